$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q5) updated values
$ws.Range("B7").Value = -0.02926116508521665
$ws.Range("C7").Value = 0.4237601240853551
$ws.Range("D7").Value = 0.2924717721126746
$ws.Range("E7").Value = 0.5408065939988848
$ws.Range("F7").Value = 0.5727717703168048
$ws.Range("G7").Value = 9

# Row 8 (Q6) updated values
$ws.Range("B8").Value = -0.0867339457287771
$ws.Range("C8").Value = 0.340876242626539
$ws.Range("D8").Value = 0.1524713185586605
$ws.Range("E8").Value = 0.3904757592458979
$ws.Range("F8").Value = 0.4038156867546136
$ws.Range("G8").Value = 9
